# Update cryptos list figures (prices/volumes) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.924.24"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.515.23"
$ws.Range("D5").Value = "'588.28"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'133.67"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "3.513.34"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "4.113.85"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'27.58"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D17").Value = "3.516.82"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "64.931.70"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'10.02"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "'5.64"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "'390.40"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'74.90"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "3.657.44"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("E28").Value = "  +8.12%  "
$ws.Range("D29").Value = "'7.56"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "3.521.43"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'169.84"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.57"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").Value = "'6.93"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "'26.09"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").Value = "'42.90"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Value = "2.489.31"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("E51").Value = "  +3.56%  "
